$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Long repr text of the model, reused for rows 2-5.
$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n" +
    "                                            estimator=Pipeline(steps=[('model',`n" +
    "                                                                       AdaBoostRegressor())]),`n" +
    "                                            param_grid={'model__learning_rate': [0.1,`n" +
    "                                                                                 0.5,`n" +
    "                                                                                 1.0],`n" +
    "                                                        'model__n_estimators': [50,`n" +
    "                                                                                100,`n" +
    "                                                                                150]},`n" +
    "                                            scoring='neg_mean_squared_error'))"

# New header cell F1: reuse the existing header style (bold/border/center)
# by copying formats from A1, then set its text.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Updated numeric metrics for rows 2-5 (columns B, C, D).
$ws.Range("B2").Value = 0.6648436607028957
$ws.Range("C2").Value = 0.9804718000845237
$ws.Range("D2").Value = 0.5985497708647042

$ws.Range("B3").Value = 0.4007794473230393
$ws.Range("C3").Value = 0.994461119521201
$ws.Range("D3").Value = 0.5312248764015917

$ws.Range("B4").Value = 0.1553632948569371
$ws.Range("C4").Value = 0.9983987036474841
$ws.Range("D4").Value = 0.3322693132296734

$ws.Range("B5").Value = 0.5006861154425049
$ws.Range("C5").Value = 0.9970175556674374
$ws.Range("D5").Value = 0.5868457625759361

# New column F with model description for rows 2-5.
$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText
$ws.Range("F4").Value = $modelText
$ws.Range("F5").Value = $modelText

# The engine auto-sizes row height whenever a cell with embedded newlines is
# written; re-run AutoFit so the rows fall back to the sheet default height
# instead of keeping a pinned / explicit height.
$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(3).EntireRow.AutoFit()
$ws.Rows.Item(4).EntireRow.AutoFit()
$ws.Rows.Item(5).EntireRow.AutoFit()
